$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.227987442866985
$ws.Range("C2").Value = 0.908388657222458
$ws.Range("D2").Value = 1.06954675481358
$ws.Range("E2").Value = 0.910832714319771
$ws.Range("F2").Value = 0.794800886664094
$ws.Range("G2").Value = 0.851526975572939

$ws.Range("B3").Value = 1.29012170421717
$ws.Range("C3").Value = 0.71373160713859
$ws.Range("D3").Value = 0.496042772941757
$ws.Range("E3").Value = 0.799851834930106
$ws.Range("F3").Value = 0.834655917235992
$ws.Range("G3").Value = 0.771063245736172

$ws.Range("B4").Value = 1.10926326827653
$ws.Range("C4").Value = 0.51466722089331
$ws.Range("D4").Value = 0.211105836973105
$ws.Range("E4").Value = 0.74432353515363
$ws.Range("F4").Value = 0.723279278626472
$ws.Range("G4").Value = 0.576496484207996

$ws.Range("B5").Value = 0.654512957312685
$ws.Range("C5").Value = 0.798072488788741
$ws.Range("D5").Value = 0.968025998751942
$ws.Range("E5").Value = 0.723116256852647
$ws.Range("F5").Value = 0.817536514958359
$ws.Range("G5").Value = 0.70834869538422

$ws.Range("B6").Value = 0.379813210645749
$ws.Range("C6").Value = 0.859235379980694
$ws.Range("D6").Value = 1.11214525010227
$ws.Range("E6").Value = 0.69084106402753
$ws.Range("F6").Value = 0.870990343302749
$ws.Range("G6").Value = 0.780455656279622

$ws.Range("B7").Value = 0.524413356901235
$ws.Range("C7").Value = 0.585625605236319
$ws.Range("D7").Value = 0.728948201773656
$ws.Range("E7").Value = 0.383382498494704
$ws.Range("F7").Value = 1.05151585827848
$ws.Range("G7").Value = 0.52659867237478

$ws.Range("B8").Value = 0.514200814604642
$ws.Range("C8").Value = 0.779933020673987
$ws.Range("D8").Value = 0.823717709544782
$ws.Range("E8").Value = 0.709265859816006
$ws.Range("F8").Value = 0.93289574312254
$ws.Range("G8").Value = 0.734479715096122

$ws.Range("B9").Value = 0.855893363074046
$ws.Range("C9").Value = 0.797391506851663
$ws.Range("D9").Value = 0.980991011502419
$ws.Range("E9").Value = 0.671367010266482
$ws.Range("F9").Value = 1.08020704088589
$ws.Range("G9").Value = 0.581903795122262

$ws.Range("B10").Value = 1.19208483609783
$ws.Range("C10").Value = 0.949770482384599
$ws.Range("D10").Value = 0.9875534552433
$ws.Range("E10").Value = 0.950191641790629
$ws.Range("F10").Value = 1.10893708753354
$ws.Range("G10").Value = 0.839781045872143

$ws.Range("B11").Value = 1.22411264136072
$ws.Range("C11").Value = 0.69929992994411
$ws.Range("D11").Value = 0.735244891461985
$ws.Range("E11").Value = 0.580082836498676
$ws.Range("F11").Value = 1.04728151255664
$ws.Range("G11").Value = 0.614426923044058

$ws.Range("B12").Value = 1.33752926205328
$ws.Range("C12").Value = 0.999997430366336
$ws.Range("D12").Value = 1.16673423948734
$ws.Range("E12").Value = 0.54040163604503
$ws.Range("F12").Value = 1.34326922607862
$ws.Range("G12").Value = 1.14700450104329

